# Applies a cyclic rotation of the species-observation rows (2 -> 6 -> 4 -> 3 -> 2)
# in the "Artfynd" worksheet. Only the observation-specific columns are moved;
# the shared/common columns (C, D, I, K, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AG, AT, AW, AX, AY) are identical across rows 2, 3, 4 and 6 and therefore
# need no change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually differ between the affected rows.
$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R")

# Capture the "before" values for rows 2, 3, 4 and 6 so we can redistribute them.
$rows = @(2, 3, 4, 6)
$values = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $values[$r] = $rowVals
}

# New row 2 gets old row 3's data
# New row 3 gets old row 4's data
# New row 4 gets old row 6's data
# New row 6 gets old row 2's data
$mapping = @{
    2 = 3
    3 = 4
    4 = 6
    6 = 2
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $values[$srcRow][$col]
    }
}
